$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("N2").Value = 4.1
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.83
$ws.Range("R2").Value = 1.42
$ws.Range("S2").Value = 3.1
$ws.Range("T2").Value = 1.82
$ws.Range("U2").Value = 2.06
$ws.Range("AB2").Value = 9.199999999999999
$ws.Range("AE2").Value = 90
$ws.Range("AN2").Value = 9.800000000000001

# Row 4 updates
$ws.Range("N4").Value = 3.8
